$wb = $excel.ActiveWorkbook

function Set-Cells {
    param($ws, $row, $values)
    foreach ($key in $values.Keys) {
        $ws.Range("$key$row").Value = $values[$key]
    }
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-Cells $ws 34 @{ "H" = 809.8461; "I" = 78.666664; "J" = 905.2174; "K" = 235.999992; "L" = 2715.6522; "M" = -151.999992; "N" = -2883.6522 }
Set-Cells $ws 36 @{ "H" = 2639.8; "J" = 3174.75; "L" = 9524.25; "N" = -9862.25 }
Set-Cells $ws 51 @{ "H" = 1353.4286; "I" = 400; "J" = 1734.8; "K" = 1200; "L" = 5204.4; "M" = -740; "N" = -6124.4 }
Set-Cells $ws 120 @{ "H" = 12060; "I" = 7433.3335; "J" = 19000; "K" = 22300.0005; "L" = 57000; "M" = -17462.0005; "N" = -66676 }
Set-Cells $ws 121 @{ "H" = 58706.543; "I" = 10176.667; "J" = 63256.22; "K" = 30530.001; "L" = 189768.66; "M" = -29220.001; "N" = -192388.66 }
Set-Cells $ws 122 @{ "H" = 329.35898; "I" = 271.1613; "J" = 554.875; "K" = 2440.4517; "L" = 4993.875; "M" = 9.548300000000381; "N" = -9893.875 }
Set-Cells $ws 123 @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }
Set-Cells $ws 124 @{ "H" = 1864.5454; "I" = 0; "J" = 1864.5454; "K" = 0; "L" = 5593.6362; "N" = -15413.6362 }
Set-Cells $ws 125 @{ "H" = 6000; "I" = 0; "J" = 6000; "K" = 0; "L" = 18000; "N" = -27840 }
Set-Cells $ws 126 @{ "H" = 1773.3334; "I" = 910; "J" = 3500; "K" = 2730; "L" = 10500; "M" = 2210; "N" = -20380 }
Set-Cells $ws 127 @{ "H" = 1348.5; "I" = 1300; "J" = 1353.8889; "K" = 3900; "L" = 4061.6667; "M" = 1060; "N" = -13981.6667 }
Set-Cells $ws 128 @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }
Set-Cells $ws 129 @{ "H" = 14493983; "I" = 1018.3333; "J" = 19609146; "K" = 3054.9999; "L" = 58827438; "M" = 1945.0001; "N" = -58837438 }
Set-Cells $ws 130 @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }
Set-Cells $ws 131 @{ "H" = 6025043; "I" = 300.72726; "J" = 6945490; "K" = 902.18178; "L" = 20836470; "M" = 4137.81822; "N" = -20846550 }
Set-Cells $ws 132 @{ "H" = 976.64703; "I" = 638; "J" = 1161.3636; "K" = 5742; "L" = 10452.2724; "M" = -3212; "N" = -15512.2724 }
Set-Cells $ws 133 @{ "H" = 66672300; "I" = 90914060; "J" = 7450; "K" = 272742180; "L" = 22350; "M" = -272737120; "N" = -32470 }
Set-Cells $ws 134 @{ "H" = 2826.25; "I" = 2230; "J" = 7000; "K" = 6690; "L" = 21000; "M" = -1620; "N" = -31140 }
Set-Cells $ws 136 @{ "H" = 3000; "I" = 3000; "J" = 0; "K" = 9000; "L" = 0; "M" = -3900 }
Set-Cells $ws 137 @{ "H" = 8814.235000000001; "I" = 15260; "J" = 4302.2; "K" = 45780; "L" = 12906.6; "M" = -40680; "N" = -23106.6 }
Set-Cells $ws 138 @{ "H" = 1252.6; "I" = 557.5; "J" = 4033; "K" = 1672.5; "L" = 12099; "M" = 3467.5; "N" = -22379 }
Set-Cells $ws 139 @{ "H" = 2178; "I" = 1722.5; "J" = 4000; "K" = 5167.5; "L" = 12000; "M" = -27.5; "N" = -22280 }
Set-Cells $ws 140 @{ "H" = 6357.263; "I" = 5385.933; "J" = 9999.75; "K" = 16157.799; "L" = 29999.25; "M" = -10977.799; "N" = -40359.25 }
Set-Cells $ws 141 @{ "H" = 7463; "I" = 8518.571; "J" = 5000; "K" = 25555.713; "L" = 15000; "M" = -20375.713; "N" = -25360 }

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

Set-Cells $ws 119 @{ "H" = 30870; "I" = 0; "J" = 30870; "K" = 0; "L" = 30870; "N" = -40546 }
Set-Cells $ws 120 @{ "H" = 27000; "I" = 0; "J" = 27000; "K" = 0; "L" = 27000; "N" = -36676 }
Set-Cells $ws 121 @{ "H" = 19000; "I" = 19000; "J" = 0; "K" = 19000; "L" = 0; "M" = -17253 }
Set-Cells $ws 122 @{ "H" = 1701.8; "I" = 2004; "J" = 1626.25; "K" = 6012; "L" = 4878.75; "M" = -3562; "N" = -9778.75 }
Set-Cells $ws 123 @{ "H" = 37124.348; "I" = 0; "J" = 37124.348; "K" = 0; "L" = 37124.348; "N" = -46924.348 }
Set-Cells $ws 124 @{ "H" = 27164.25; "I" = 0; "J" = 27164.25; "K" = 0; "L" = 27164.25; "N" = -36984.25 }
Set-Cells $ws 125 @{ "H" = 29950; "I" = 0; "J" = 29950; "K" = 0; "L" = 29950; "N" = -39790 }
Set-Cells $ws 126 @{ "H" = 2572.5715; "I" = 2002; "J" = 3333.3333; "K" = 6006; "L" = 9999.999899999999; "M" = -3536; "N" = -14939.9999 }
Set-Cells $ws 127 @{ "H" = 57469.668; "I" = 0; "J" = 57469.668; "K" = 0; "L" = 57469.668; "N" = -67389.66800000001 }
Set-Cells $ws 128 @{ "H" = 0; "I" = 0; "J" = 0; "K" = 0; "L" = 0 }
Set-Cells $ws 129 @{ "H" = 30000; "I" = 0; "J" = 30000; "K" = 0; "L" = 30000; "N" = -40000 }
Set-Cells $ws 130 @{ "H" = 48999.668; "I" = 0; "J" = 48999.668; "K" = 0; "L" = 48999.668; "N" = -59039.668 }
Set-Cells $ws 131 @{ "H" = 44990; "I" = 0; "J" = 44990; "K" = 0; "L" = 44990; "N" = -55070 }
Set-Cells $ws 132 @{ "H" = 56252224; "I" = 93751240; "J" = 3697.0625; "K" = 281253720; "L" = 11091.1875; "M" = -281251190; "N" = -16151.1875 }
Set-Cells $ws 133 @{ "H" = 42107.5; "I" = 0; "J" = 42107.5; "K" = 0; "L" = 42107.5; "N" = -52227.5 }
Set-Cells $ws 135 @{ "H" = 55360.824; "I" = 0; "J" = 55360.824; "K" = 0; "L" = 55360.824; "N" = -65500.824 }
Set-Cells $ws 136 @{ "H" = 77321.30499999999; "I" = 100372.2; "J" = 485; "K" = 301116.6; "L" = 1455; "M" = -298566.6; "N" = -6555 }
Set-Cells $ws 137 @{ "H" = 50000; "I" = 0; "J" = 50000; "K" = 0; "L" = 50000; "N" = -60200 }
Set-Cells $ws 138 @{ "H" = 42900; "I" = 0; "J" = 42900; "K" = 0; "L" = 42900; "N" = -53180 }
Set-Cells $ws 139 @{ "H" = 48650; "I" = 0; "J" = 48650; "K" = 0; "L" = 48650; "N" = -58930 }
Set-Cells $ws 140 @{ "H" = 53266.668; "I" = 0; "J" = 53266.668; "K" = 0; "L" = 53266.668; "N" = -63626.668 }
Set-Cells $ws 141 @{ "H" = 49996.668; "I" = 0; "J" = 49996.668; "K" = 0; "L" = 49996.668; "N" = -60356.668 }
